$wb = $excel.ActiveWorkbook

# Target F-column (想去人数 / "want to go" count) values to apply
# to both the "展览" (sheet 1) and "全部类型" (sheet 4) sheets.
$targetValues = @{
    3 = 3160
    4 = 229
    5 = 123
    7 = 1671
    9 = 461
    10 = 367
    14 = 32
    15 = 229
    16 = 237
    17 = 230
    18 = 6
    19 = 23
    21 = 50
    23 = 374
    24 = 191
    26 = 28
    27 = 7
    28 = 21
    29 = 168
    30 = 2129
    33 = 464
    34 = 214
    36 = 424
    37 = 225
    39 = 412
    40 = 509
    41 = 414
}

$sheetNames = @("展览", "全部类型")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $targetValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $targetValues[$row]
    }
}

Write-Output "Done updating F-column want-to-go counts."
